# Updated cryptos list - refresh Price (col D) and Volume(1h) (col E) values,
# plus the Monero/EthereumClassic row-content swap (rows 27-28).
# Column D values that look numeric are forced to Text ("@") first so that
# exact-looking-numeric strings (with meaningful trailing zeros, e.g.
# "20.90", "0.07680") are preserved verbatim instead of being normalized
# into a floating point number by Excel's automatic type detection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.429.57'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.820.77'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.49'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5419'
$ws.Range("E7").Value = '  +1.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4032'
$ws.Range("E8").Value = '  +7.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07680'
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.117'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.86'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.319'
$ws.Range("E12").Value = '  +3.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.638'
$ws.Range("E13").Value = '  +5.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.002'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.90'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '1.826.15'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("E17").Value = '  +2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.66'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06603'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.65'
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.060'
$ws.Range("E22").Value = '  +2.66%  '
$ws.Range("D23").Value = '28.438.15'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.273'
$ws.Range("E25").Value = '  +8.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.463'
$ws.Range("E26").Value = '  +7.91%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.74'
$ws.Range("E27").Value = '  +2.55%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.39'
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").Value = '2.037.95'
$ws.Range("E29").Value = '  +2.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.59'
$ws.Range("E30").Value = '  +2.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1112'
$ws.Range("E31").Value = '  +5.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.121'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.678'
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07376'
$ws.Range("E34").Value = '  +12.99%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2241'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02336'
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.201'
$ws.Range("E38").Value = '  +3.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.837'
$ws.Range("E39").Value = '  +4.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.34'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6280'
$ws.Range("E41").Value = '  +1.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.179'
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.400'
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.48'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5849'
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.87'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.003'
$ws.Range("E49").Value = '  +3.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.197'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06870'
$ws.Range("E51").Value = '  +0.81%  '
